$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: tc14 continues
$ws.Range("A31").Value = 41466
$ws.Range("B31").Value = 2.5
$ws.Range("C31").Value = 1.5
$ws.Range("D31").Value = "Implementation tc14"

# Row 32: tc14 continues
$ws.Range("A32").Value = 41467
$ws.Range("C32").Value = 4.25
$ws.Range("D32").Value = "Implementation tc14"

# Carry the date number-format down into the newly added date cells
$ws.Range("A30").Copy() | Out-Null
$ws.Range("A31:A32").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C32").Select() | Out-Null
